$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.019999999999999
$ws.Range("C2").Value = 1.045888919939243
$ws.Range("D2").Value = 1.051881940192187
$ws.Range("E2").Value = 1.053268210991737
$ws.Range("F2").Value = 1.063311325738323
$ws.Range("I2").Value = 1.044905204617403
$ws.Range("J2").Value = 1.050946246141701
$ws.Range("K2").Value = 1.054632212339488
$ws.Range("L2").Value = 1.056014648806957
$ws.Range("M2").Value = 1.066030308212694
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.046711308515427
$ws.Range("D3").Value = 1.052508943949808
$ws.Range("E3").Value = 1.053979968179239
$ws.Range("F3").Value = 1.064056309778269
$ws.Range("I3").Value = 1.045091769524372
$ws.Range("J3").Value = 1.051417107097216
$ws.Range("K3").Value = 1.055072428311692
$ws.Range("L3").Value = 1.056539676370491
$ws.Range("M3").Value = 1.066590454118423
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.047244254333979
$ws.Range("D4").Value = 1.052915360384767
$ws.Range("E4").Value = 1.05444158529646
$ws.Range("F4").Value = 1.06453943265089
$ws.Range("I4").Value = 1.045211735609935
$ws.Range("J4").Value = 1.051721900334248
$ws.Range("K4").Value = 1.055357287291028
$ws.Range("L4").Value = 1.056879790596831
$ws.Range("M4").Value = 1.066953303789797
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.047468495697132
$ws.Range("D5").Value = 1.053086384250074
$ws.Range("E5").Value = 1.054635901507857
$ws.Range("F5").Value = 1.064742791170015
$ws.Range("I5").Value = 1.045261988085564
$ws.Range("J5").Value = 1.051850061463763
$ws.Range("K5").Value = 1.055477042945489
$ws.Range("L5").Value = 1.057022865531827
$ws.Range("M5").Value = 1.067105939343361
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.047506157951058
$ws.Range("D6").Value = 1.053115109611215
$ws.Range("E6").Value = 1.054668542795457
$ws.Range("F6").Value = 1.064776950773774
$ws.Range("I6").Value = 1.045270415044979
$ws.Range("J6").Value = 1.051871581768096
$ws.Range("K6").Value = 1.055497150460956
$ws.Range("L6").Value = 1.057046893728763
$ws.Range("M6").Value = 1.06713157296132
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.04724724991013
$ws.Range("D7").Value = 1.052917644962622
$ws.Range("E7").Value = 1.054444180770066
$ws.Range("F7").Value = 1.064542148943306
$ws.Range("I7").Value = 1.045212407799208
$ws.Range("J7").Value = 1.051723612728933
$ws.Range("K7").Value = 1.055358887470021
$ws.Range("L7").Value = 1.056881702015849
$ws.Range("M7").Value = 1.066955342948058
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.046166682111975
$ws.Range("D8").Value = 1.052093692315917
$ws.Range("E8").Value = 1.053508531429313
$ws.Range("F8").Value = 1.063562874098921
$ws.Range("I8").Value = 1.04496841072533
$ws.Range("J8").Value = 1.05110535120631
$ws.Range("K8").Value = 1.054780982510441
$ws.Range("L8").Value = 1.056192003407663
$ws.Range("M8").Value = 1.06621952876059
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.044268835320233
$ws.Range("D9").Value = 1.05064725219643
$ws.Range("E9").Value = 1.051868025452025
$ws.Range("F9").Value = 1.061845542554791
$ws.Range("I9").Value = 1.044532717664341
$ws.Range("J9").Value = 1.050016836369523
$ws.Range("K9").Value = 1.053762777592207
$ws.Range("L9").Value = 1.054979693975409
$ws.Range("M9").Value = 1.064926051517513
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.043007919792841
$ws.Range("D10").Value = 1.049686751307041
$ws.Range("E10").Value = 1.050780007453413
$ws.Range("F10").Value = 1.060706345338236
$ws.Range("I10").Value = 1.04423844411677
$ws.Range("J10").Value = 1.049291876821326
$ws.Range("K10").Value = 1.053084151040741
$ws.Range("L10").Value = 1.054173616833449
$ws.Range("M10").Value = 1.064065932882405
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.042462976828691
$ws.Range("D11").Value = 1.04927176775148
$ws.Range("E11").Value = 1.050310250015834
$ws.Range("F11").Value = 1.060214436864409
$ws.Range("I11").Value = 1.044110126356056
$ws.Range("J11").Value = 1.048978149267948
$ws.Range("K11").Value = 1.052790358035928
$ws.Range("L11").Value = 1.053825101364986
$ws.Range("M11").Value = 1.063694035108966
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.042260719187369
$ws.Range("D12").Value = 1.049117764540866
$ws.Range("E12").Value = 1.050135967836318
$ws.Range("F12").Value = 1.060031928385335
$ws.Range("I12").Value = 1.044062329703447
$ws.Range("J12").Value = 1.048861646084572
$ws.Range("K12").Value = 1.05268124020926
$ws.Range("L12").Value = 1.053695727126297
$ws.Range("M12").Value = 1.063555978475601
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.042304096937394
$ws.Range("D13").Value = 1.049150792377916
$ws.Range("E13").Value = 1.050173342570563
$ws.Range("F13").Value = 1.060071067615682
$ws.Range("I13").Value = 1.044072588286116
$ws.Range("J13").Value = 1.048886635068348
$ws.Range("K13").Value = 1.0527046458788
$ws.Range("L13").Value = 1.053723474697077
$ws.Range("M13").Value = 1.0635855883192
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.042446254905719
$ws.Range("D14").Value = 1.049259034927353
$ws.Range("E14").Value = 1.050295839568838
$ws.Range("F14").Value = 1.060199346396119
$ws.Range("I14").Value = 1.044106178193134
$ws.Range("J14").Value = 1.048968518474452
$ws.Range("K14").Value = 1.05278133811045
$ws.Range("L14").Value = 1.053814405610318
$ws.Range("M14").Value = 1.063682621606285
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.04253386417412
$ws.Range("D15").Value = 1.049325745366405
$ws.Range("E15").Value = 1.050371341457922
$ws.Range("F15").Value = 1.060278410844399
$ws.Range("I15").Value = 1.04412685634773
$ws.Range("J15").Value = 1.049018973455844
$ws.Range("K15").Value = 1.052828592090794
$ws.Range("L15").Value = 1.053870441785946
$ws.Range("M15").Value = 1.063742418035253
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.043044107975871
$ws.Range("D16").Value = 1.049714311945132
$ws.Range("E16").Value = 1.050811212572181
$ws.Range("F16").Value = 1.060739020764968
$ws.Range("I16").Value = 1.044246941347435
$ws.Range("J16").Value = 1.049312701901086
$ws.Range("K16").Value = 1.053103650434528
$ws.Range("L16").Value = 1.054196757777076
$ws.Range("M16").Value = 1.064090626023109
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.04336445070323
$ws.Range("D17").Value = 1.049958297120638
$ws.Range("E17").Value = 1.051087498101284
$ws.Range("F17").Value = 1.06102831791676
$ws.Range("I17").Value = 1.044322028340298
$ws.Range("J17").Value = 1.049497000469779
$ws.Range("K17").Value = 1.053276203487493
$ws.Range("L17").Value = 1.054401587869434
$ws.Range("M17").Value = 1.064309193300606
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.043551401583834
$ws.Range("D18").Value = 1.050100698146964
$ws.Range("E18").Value = 1.051248782004224
$ws.Range("F18").Value = 1.061197192228252
$ws.Range("I18").Value = 1.044365738868854
$ws.Range("J18").Value = 1.049604516432074
$ws.Range("K18").Value = 1.053376856144381
$ws.Range("L18").Value = 1.054521111982223
$ws.Range("M18").Value = 1.064436731811018
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.043615163971864
$ws.Range("D19").Value = 1.050149268184732
$ws.Range("E19").Value = 1.051303797842028
$ws.Range("F19").Value = 1.061254796357368
$ws.Range("I19").Value = 1.044380628342895
$ws.Range("J19").Value = 1.049641179559233
$ws.Range("K19").Value = 1.053411176967862
$ws.Range("L19").Value = 1.05456187507917
$ws.Range("M19").Value = 1.064480227899994
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.043330070562861
$ws.Range("D20").Value = 1.04993211064156
$ws.Range("E20").Value = 1.05105784169011
$ws.Range("F20").Value = 1.060997265382759
$ws.Range("I20").Value = 1.04431398115215
$ws.Range("J20").Value = 1.049477225133701
$ws.Range("K20").Value = 1.053257689610974
$ws.Range("L20").Value = 1.05437960633649
$ws.Range("M20").Value = 1.064285737717849
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.042404388538084
$ws.Range("D21").Value = 1.049227156309381
$ws.Range("E21").Value = 1.050259761534689
$ws.Range("F21").Value = 1.060161565726816
$ws.Range("I21").Value = 1.044096290487725
$ws.Range("J21").Value = 1.048944405034597
$ws.Range("K21").Value = 1.052758753872714
$ws.Range("L21").Value = 1.053787626503534
$ws.Range("M21").Value = 1.063654045423233
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.041823293143252
$ws.Range("D22").Value = 1.04878473544088
$ws.Range("E22").Value = 1.049759173505505
$ws.Range("F22").Value = 1.059637334021669
$ws.Range("I22").Value = 1.043958646442871
$ws.Range("J22").Value = 1.048609569541535
$ws.Range("K22").Value = 1.052445111608127
$ws.Range("L22").Value = 1.053415888499151
$ws.Range("M22").Value = 1.063257355080561
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.042131255182536
$ws.Range("D23").Value = 1.04901919349194
$ws.Range("E23").Value = 1.050024430565222
$ws.Range("F23").Value = 1.059915124138204
$ws.Range("I23").Value = 1.044031687187259
$ws.Range("J23").Value = 1.048787055607519
$ws.Range("K23").Value = 1.052611373284596
$ws.Range("L23").Value = 1.053612909403804
$ws.Range("M23").Value = 1.063467602088848
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.043345605158487
$ws.Range("D24").Value = 1.049943942910777
$ws.Range("E24").Value = 1.051071241740585
$ws.Range("F24").Value = 1.061011296278124
$ws.Range("I24").Value = 1.044317617597289
$ws.Range("J24").Value = 1.049486160702279
$ws.Range("K24").Value = 1.05326605521842
$ws.Range("L24").Value = 1.054389538690194
$ws.Range("M24").Value = 1.064296336125886
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.044758721779053
$ws.Range("D25").Value = 1.051020531590414
$ws.Range("E25").Value = 1.052291148088256
$ws.Range("F25").Value = 1.062288520174852
$ws.Range("I25").Value = 1.044646030023766
$ws.Range("J25").Value = 1.050298123345534
$ws.Range("K25").Value = 1.054025982992344
$ws.Range("L25").Value = 1.055292736591284
$ws.Range("M25").Value = 1.065260066091853
